$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Medicine" column (E) ------------------------------------------------

# Header cell E1: same green fill as the other headers (A1:D1) but using the
# "Arial (Body)" font family that the sheet author picked for the new column.
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E1").Font.Name = "Arial (Body)"
$ws.Range("E1").Value = "Medicine"

# Body cells E2:E13: reuse the wrap-text body style already used in column D
# (D13) so we don't invent a new look for the medicine list.
$ws.Range("D13").Copy()

$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E2").Value = "Topical retinoids`nTopical antibiotics`nIsotretinoin`nBenzoyl Peroxide Cream`nSalicylic Acid Cream"

$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = "Topical corticosteroids`nBarrier creams`nEmollients`nOral Antihistamines`nSystemic Immunosuppressants"

$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E4").Value = "Biologic Drugs`nMethotrexate`nCyclosporine`nClobetasol`nCalcipotriene`nCoal Tar"

$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Value = "Imiquimod`n5-Fluorouracil`nPhotodynamic Therapy`nVismodegib`nSonidegib"

$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Value = "Chemotherapy`nImmunotherapy`nEfudex`nAldara`nZyclara"

$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E7").Value = "Topical corticosteroid`nTopical emollients`nImmunotherapy`nChemotherapy"

# Rows 9-13 were typed in before row 8 got its medicine note, so fill them in
# this order to line up with the shared-string table the author ended up with.
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").Value = "antihistamines`ncetirizine`ndiphenhydramie`nBenadryl"

$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").Value = "Griseofulvin (Grifulvin V, Gris-PEG)`nTerbinafine.`nItraconazole (Onmel, Sporanox)`nFluconazole (Diflucan)"

$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("E11").Value = "iodine `nsalicylic acid`npotassium hydroxide`ntretinoin cantharidin"

$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("E12").Value = "Eskata"

$ws.Range("E13").PasteSpecial(-4122)
$ws.Range("E13").Value = "efinaconazole`ntavaborole"

$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E8").Value = "Melanocytic nevi can be surgically removed for cosmetic considerations or because of concern regarding the biological potential of a lesion"

# Column E width, matching the width used for the rest of the medicine column.
$ws.Columns.Item(5).ColumnWidth = 73.45

# Put the selection/view on the new column, like the author did after adding it.
$ws.Range("E1").Select()
